$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.460.74'
$ws.Range("E2").Value = '  -0.67%  '
$ws.Range("D3").Value = '3.057.65'
$ws.Range("E3").Value = '  -2.75%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.25'
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.68'
$ws.Range("E6").Value = '  +6.35%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.547'
$ws.Range("E8").Value = '  +3.35%  '
$ws.Range("D9").Value = '3.075.37'
$ws.Range("E9").Value = '  -1.97%  '
$ws.Range("E10").Value = '  -1.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.85'
$ws.Range("E11").Value = '  -0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.461'
$ws.Range("E12").Value = '  -0.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.61'
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("E14").Value = '  -2.07%  '
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("D16").Value = '3.564.35'
$ws.Range("E16").Value = '  -2.80%  '
$ws.Range("E17").Value = '  -1.61%  '
$ws.Range("D18").Value = '63.406.35'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("D19").Value = '3.069.63'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '477.30'
$ws.Range("E20").Value = '  +1.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.74'
$ws.Range("E21").Value = '  +2.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.719'
$ws.Range("E22").Value = '  -1.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.57'
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("E24").Value = '  +2.85%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.94'
$ws.Range("E25").Value = '  -0.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.29'
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.04'
$ws.Range("E27").Value = '  +1.64%  '
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("E32").Value = '  -2.08%  '
$ws.Range("E33").Value = '  +2.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.27'
$ws.Range("E34").Value = '  -1.74%  '
$ws.Range("D35").Value = '0.0₃0849'
$ws.Range("E35").Value = '  +0.92%  '
$ws.Range("E36").Value = '  -1.69%  '
$ws.Range("E37").Value = '  +4.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.13'
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("E39").Value = '  -3.85%  '
$ws.Range("E40").Value = '  +1.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.59'
$ws.Range("E41").Value = '  -1.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '446.30'
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.286'
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("E44").Value = '  -2.15%  '
$ws.Range("B45").Value = 'Kaspa'
$ws.Range("C45").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.112'
$ws.Range("E45").Value = '  +3.56%  '
$ws.Range("B46").Value = 'Arweave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.12'
$ws.Range("E46").Value = '  +1.36%  '
$ws.Range("D47").Value = '2.807.72'
$ws.Range("E47").Value = '  -3.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '132.21'
$ws.Range("E48").Value = '  +1.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '25.23'
$ws.Range("E50").Value = '  +3.94%  '
$ws.Range("E51").Value = '  +0.75%  '
